$wb = $excel.ActiveWorkbook

# ---- Sheet "log sheet" (sheet1) ----
$ws1 = $wb.Worksheets.Item("log sheet")

# New row 10: task "Homepage Template Integration" completed 28/09/2013, 4 hrs @13 = 52 USD
$ws1.Range("A10").Value = 5
$ws1.Range("B10").Value = "Homepage Template Integration"
$ws1.Range("C10").Value = "myguitarpal"
$ws1.Range("D10").Value = "28/09/2013"
$ws1.Range("E10").Value = 4
$ws1.Range("F10").Value = 13
$ws1.Range("G10").Value = 52
$ws1.Range("H10").Value = "USD"

# ---- Sheet "payment" (sheet2) ----
$ws2 = $wb.Worksheets.Item("payment")

# Existing Safepay Guru payment (18/09/2013) confirmed released, 195 USD, no action pending now
$ws2.Range("B9").Value = 195
$ws2.Range("C9").Value = "Released"
$ws2.Range("F9").Value = "NULL"

# New payment received on Safepay Guru: 200 USD on 29/09/2013, funded, release action still required
$ws2.Range("A10").Value = "29/09/2013"
$ws2.Range("B10").Value = 200
$ws2.Range("C10").Value = "Funded"
$ws2.Range("D10").Value = "USD"
$ws2.Range("E10").Value = "Safepay Guru"
$ws2.Range("F10").Value = "Release"

# Total payment received updated
$ws2.Range("B16").Value = 195

# Title changed to the new project (set last so the new string appends at the
# end of the shared-string table, matching the order they were introduced)
$ws1.Range("A1").Value = "Project: myguitarpal"

# Undo the auto row-height recalculation that setting the big banner cell's
# value triggers, so row 1 keeps its original explicit height
$ws1.Rows.Item(1).RowHeight = 15

# Update the selection/view state: payment sheet lands on B16, then focus
# returns to the log sheet with the header block selected
[void]$ws2.Range("B16").Select()
[void]$ws1.Range("A1:N4").Select()
